# "modification of balance titles"
#
# SoldeArmateur.xlsx has two sheets:
#   1) "Armateurs Créditeurs" - list of armateurs with a credit balance
#   2) "Armateurs Débiteurs"  - list of armateurs with a debit balance
#
# The edit:
#  - refreshes the "Édité le : <date> à <time>" export stamp (shared by
#    both sheets' A1 cell)
#  - retitles sheet 1's banner from "Soldes Des Armateurs" to
#    "Solde des Armateurs Créditeurs"
#  - retitles sheet 2's banner to the distinct "Soldes des Armateurs
#    Débiteurs" (previously it reused sheet 1's banner text)
#  - turns the two balance amounts on sheet 1 ("653,00" / "39 141,00",
#    stored as text) into real numbers (653 / 39141) formatted with
#    #,##0.00 and right-aligned, matching how the "Total" row already
#    renders its amount

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$newStamp = "Édité le : 28/01/2025 à 09:16:13 `n par :"

# --- Sheet 1: "Armateurs Créditeurs" ---
$ws1.Range("A1").Value = $newStamp
$ws1.Range("A2").Value = "Solde des Armateurs Créditeurs"

$ws1.Range("B6").Value = 653
$ws1.Range("B6").NumberFormat = "#,##0.00"
$ws1.Range("B6").HorizontalAlignment = -4152

$ws1.Range("B7").Value = 39141
$ws1.Range("B7").NumberFormat = "#,##0.00"
$ws1.Range("B7").HorizontalAlignment = -4152

# --- Sheet 2: "Armateurs Débiteurs" ---
$ws2.Range("A1").Value = $newStamp
$ws2.Range("A2").Value = "Soldes des Armateurs Débiteurs"
